$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Gantt chart dates: merge split "Month" + " " + "Day" runs into a single
#    run reading "Month Day" (Find/Replace naturally coalesces the matched
#    range into one run with the replacement text).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Feb 4", $true, $false, $false, $false, $false, $true, 1, $false, "Feb 4", 2)
$d.Content.Find.Execute("Feb 5", $true, $false, $false, $false, $false, $true, 1, $false, "Feb 5", 2)
$d.Content.Find.Execute("April 10", $true, $false, $false, $false, $false, $true, 1, $false, "April 10", 2)
$d.Content.Find.Execute("April 11", $true, $false, $false, $false, $false, $true, 1, $false, "April 11", 2)
$d.Content.Find.Execute("April 16", $true, $false, $false, $false, $false, $true, 1, $false, "April 16", 2)
$d.Content.Find.Execute("April 21", $true, $false, $false, $false, $false, $true, 1, $false, "April 21", 2)

# ---------------------------------------------------------------------------
# 2) Contact table (4th table in the document): Project Manager Email
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("[example@email.com]", $true, $false, $false, $false, $false, $true, 1, $false, "ibrahimshabori@gmail.com", 2)

# ---------------------------------------------------------------------------
# 3) Contact table: Phone Number -> "+1" and "0461984198" as two separate
#    runs. A plain Find/Replace (or Range.InsertAfter) coalesces same-format
#    adjacent runs into one, so the new paragraph is built explicitly via
#    InsertXML (which keeps the two <w:r> elements distinct) and the old
#    paragraph is removed afterwards.
# ---------------------------------------------------------------------------
$contactTable = $d.Tables(4)
$phoneCell = $contactTable.Cell(3, 2)
$phoneCell.Range.Text = ""

$phoneCellRange = $contactTable.Cell(3, 2).Range
$phoneCellRange.Collapse(1)
$phoneCellRange.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>+1</w:t></w:r><w:r><w:t>0461984198</w:t></w:r></w:p>')

$oldPhonePara = $contactTable.Cell(3, 2).Range.Paragraphs.Item(1)
$oldPhonePara.Range.Delete()

# ---------------------------------------------------------------------------
# 4) Contact table: Website Link (if available) Details cell was empty,
#    now reads "N/A".
# ---------------------------------------------------------------------------
$websiteCell = $d.Tables(4).Cell(5, 2)
$websiteCell.Range.Text = "N/A"
